# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# coinranking.com snapshot, as done periodically by the GitHub Actions job.
# Two pairs of adjacent coins (rows 20/21 and 40/41, plus 42/43) also
# swapped rank order, so their Coin name / Link cells are updated too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Format the Price column as Text up front so numeric-looking values
# (e.g. "1.005") are stored as strings, matching the source data which
# uses inline strings for all Price/Volume cells.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range('D2').Value = '22.469.77'
$ws.Range('E2').Value = '  -0.07%  '

# Row 3
$ws.Range('D3').Value = '1.571.98'
$ws.Range('E3').Value = '  +0.03%  '

# Row 4
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.28%  '

# Row 5
$ws.Range('D5').Value = '1.003'
$ws.Range('E5').Value = '  +0.08%  '

# Row 6
$ws.Range('D6').Value = '288.60'
$ws.Range('E6').Value = '  -0.76%  '

# Row 7
$ws.Range('D7').Value = '0.3719'
$ws.Range('E7').Value = '  +0.47%  '

# Row 8
$ws.Range('D8').Value = '48.35'
$ws.Range('E8').Value = '  -3.14%  '

# Row 9
$ws.Range('D9').Value = '0.3346'
$ws.Range('E9').Value = '  -0.90%  '

# Row 10
$ws.Range('D10').Value = '1.134'
$ws.Range('E10').Value = '  -1.30%  '

# Row 11
$ws.Range('D11').Value = '0.07490'
$ws.Range('E11').Value = '  -0.71%  '

# Row 12
$ws.Range('D12').Value = '1.005'
$ws.Range('E12').Value = '  +0.32%  '

# Row 13
$ws.Range('D13').Value = '20.94'
$ws.Range('E13').Value = '  -1.10%  '

# Row 14
$ws.Range('D14').Value = '5.979'
$ws.Range('E14').Value = '  -0.79%  '

# Row 15
$ws.Range('D15').Value = '6.936'
$ws.Range('E15').Value = '  -0.44%  '

# Row 16
$ws.Range('D16').Value = '1.579.99'
$ws.Range('E16').Value = '  +0.50%  '

# Row 17
$ws.Range('D17').Value = '0.00001117'
$ws.Range('E17').Value = '  -0.42%  '

# Row 18
$ws.Range('D18').Value = '88.45'
$ws.Range('E18').Value = '  -2.25%  '

# Row 19
$ws.Range('D19').Value = '0.06786'
$ws.Range('E19').Value = '  +0.18%  '

# Row 20
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').Value = '1.002'
$ws.Range('E20').Value = '  +0.03%  '

# Row 21
$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').Value = '6.396'
$ws.Range('E21').Value = '  +0.62%  '

# Row 22
$ws.Range('D22').Value = '16.47'
$ws.Range('E22').Value = '  +0.25%  '

# Row 23
$ws.Range('D23').Value = '12.10'
$ws.Range('E23').Value = '  -1.26%  '

# Row 24
$ws.Range('D24').Value = '22.463.25'
$ws.Range('E24').Value = '  -0.12%  '

# Row 25
$ws.Range('D25').Value = '2.401'
$ws.Range('E25').Value = '  +0.93%  '

# Row 26
$ws.Range('D26').Value = '2.573'
$ws.Range('E26').Value = '  -2.34%  '

# Row 27
$ws.Range('D27').Value = '152.67'
$ws.Range('E27').Value = '  +2.41%  '

# Row 28
$ws.Range('D28').Value = '19.79'
$ws.Range('E28').Value = '  -1.17%  '

# Row 29
$ws.Range('D29').Value = '4.998'
$ws.Range('E29').Value = '  -1.45%  '

# Row 30
$ws.Range('D30').Value = '124.55'
$ws.Range('E30').Value = '  -0.44%  '

# Row 31
$ws.Range('D31').Value = '1.755.10'
$ws.Range('E31').Value = '  +0.20%  '

# Row 32
$ws.Range('D32').Value = '1.051'
$ws.Range('E32').Value = '  -1.56%  '

# Row 33
$ws.Range('D33').Value = '6.168'
$ws.Range('E33').Value = '  -0.49%  '

# Row 34
$ws.Range('D34').Value = '2.013'
$ws.Range('E34').Value = '  -0.06%  '

# Row 35
$ws.Range('D35').Value = '9.708'
$ws.Range('E35').Value = '  -1.00%  '

# Row 36
$ws.Range('D36').Value = '0.08325'
$ws.Range('E36').Value = '  -0.26%  '

# Row 37
$ws.Range('D37').Value = '0.02456'
$ws.Range('E37').Value = '  -1.11%  '

# Row 38
$ws.Range('D38').Value = '0.2282'
$ws.Range('E38').Value = '  -0.82%  '

# Row 39
$ws.Range('D39').Value = '0.06390'
$ws.Range('E39').Value = '  -2.52%  '

# Row 40
$ws.Range('B40').Value = 'TrustWalletToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D40').Value = '1.301'
$ws.Range('E40').Value = '  -4.18%  '

# Row 41
$ws.Range('B41').Value = 'InternetComputer(DFINITY)'
$ws.Range('C41').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D41').Value = '5.395'
$ws.Range('E41').Value = '  -0.86%  '

# Row 42
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = '0.6309'
$ws.Range('E42').Value = '  +1.29%  '

# Row 43
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').Value = '11.30'
$ws.Range('E43').Value = '  +0.13%  '

# Row 44
$ws.Range('D44').Value = '1.002'
$ws.Range('E44').Value = '  +0.04%  '

# Row 45
$ws.Range('D45').Value = '13.96'
$ws.Range('E45').Value = '  -0.78%  '

# Row 46
$ws.Range('D46').Value = '0.6143'
$ws.Range('E46').Value = '  +4.82%  '

# Row 47
$ws.Range('D47').Value = '3.785'
$ws.Range('E47').Value = '  -0.58%  '

# Row 48
$ws.Range('D48').Value = '2.061'
$ws.Range('E48').Value = '  -0.64%  '

# Row 49
$ws.Range('D49').Value = '125.36'
$ws.Range('E49').Value = '  -2.79%  '

# Row 50
$ws.Range('D50').Value = '1.216'
$ws.Range('E50').Value = '  -1.66%  '

# Row 51
$ws.Range('D51').Value = '0.07272'
$ws.Range('E51').Value = '  -0.58%  '
